$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "IOB800(x2)"
$ws.Range("D9").Value = "AttachedFunctionality"
$ws.Range("D10").Value = 1
$ws.Range("N10").Value = "Other Slot Cards  (3 of 1"

$ws.Range("N10").Select()
